# Grading homework 5 (reading quiz, column J) and adding a new
# "Homework 6" column (L) with its header label.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New "Homework 6" header in column L ---
$ws.Range("L2").Value2 = "Homework 6"

# Give column L the same kind of explicit custom width the other score
# columns (F:K) already carry.
$ws.Columns.Item(12).ColumnWidth = 12.8

# --- Homework 5 ("reading quiz", column J) scores per student row ---
$ws.Range("J3").Formula  = "=13/13"
$ws.Range("J4").Value2   = 0
$ws.Range("J5").Formula  = "=13/13"
$ws.Range("J6").Formula  = "=12.5/13"
$ws.Range("J7").Formula  = "=12.5/13"
$ws.Range("J8").Formula  = "=11/13"
$ws.Range("J9").Formula  = "=13/13"
$ws.Range("J11").Formula = "=13/13"
$ws.Range("J12").Formula = "=12/13"
$ws.Range("J13").Formula = "=13/13"
$ws.Range("J14").Formula = "=13/13"
$ws.Range("J15").Formula = "=13/13"
$ws.Range("J16").Formula = "=13/13"
$ws.Range("J17").Formula = "=13/13"
$ws.Range("J19").Formula = "=13/13"

# --- Update the active selection to match the author's final cursor spot ---
$ws.Range("J10").Select() | Out-Null
